# Froze SEC1 +TE +D123
# Updates the ASR Results worksheet: column B (recognized word) and column C
# (count) for rows 2-52 are refreshed with new ASR transcription results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; B = "<their>"; C = 63 },
    @{ Row = 3; B = "<so>"; C = 63 },
    @{ Row = 4; B = "<which>"; C = 61 },
    @{ Row = 5; B = "<on>"; C = 64 },
    @{ Row = 6; B = "<whiskey>"; C = 63 },
    @{ Row = 7; B = "<may>"; C = 63 },
    @{ Row = 8; B = "<these>"; C = 62 },
    @{ Row = 9; B = "<the>"; C = 62 },
    @{ Row = 10; B = "<there>"; C = 63 },
    @{ Row = 11; B = "<about>"; C = 62 },
    @{ Row = 12; B = "<which>"; C = 60 },
    @{ Row = 13; B = "<all>"; C = 61 },
    @{ Row = 14; B = "<the>"; C = 63 },
    @{ Row = 15; B = "<they>"; C = 63 },
    @{ Row = 16; B = "<their>"; C = 62 },
    @{ Row = 17; B = "<there>"; C = 63 },
    @{ Row = 18; B = "<all>"; C = 61 },
    @{ Row = 19; B = "<the>"; C = 64 },
    @{ Row = 20; B = "<these>"; C = 63 },
    @{ Row = 21; B = "<when>"; C = 64 },
    @{ Row = 22; B = "<their>"; C = 63 },
    @{ Row = 23; B = "<them>"; C = 63 },
    @{ Row = 24; B = "<enter>"; C = 62 },
    @{ Row = 25; B = "<this>"; C = 64 },
    @{ Row = 26; B = "<what>"; C = 63 },
    @{ Row = 27; B = "<india>"; C = 64 },
    @{ Row = 28; B = "<then>"; C = 63 },
    @{ Row = 29; B = "<ald>"; C = 63 },
    @{ Row = 30; B = "<whiskey>"; C = 63 },
    @{ Row = 31; B = "<come>"; C = 63 },
    @{ Row = 32; B = "<could>"; C = 64 },
    @{ Row = 33; B = "<what>"; C = 63 },
    @{ Row = 34; B = "<like>"; C = 63 },
    @{ Row = 35; B = "<them>"; C = 63 },
    @{ Row = 36; B = "<the>"; C = 62 },
    @{ Row = 37; B = "<seven>"; C = 64 },
    @{ Row = 38; B = "<in>"; C = 63 },
    @{ Row = 39; B = "<these>"; C = 63 },
    @{ Row = 40; B = "<which>"; C = 64 },
    @{ Row = 41; B = "<word>"; C = 62 },
    @{ Row = 42; B = "<we>"; C = 64 },
    @{ Row = 43; B = "<on>"; C = 64 },
    @{ Row = 44; B = "<them>"; C = 60 },
    @{ Row = 45; B = "<was>"; C = 60 },
    @{ Row = 46; B = "<the>"; C = 63 },
    @{ Row = 47; B = "<there>"; C = 63 },
    @{ Row = 48; B = "<long>"; C = 64 },
    @{ Row = 49; B = "<there>"; C = 62 },
    @{ Row = 50; B = "<which>"; C = 64 },
    @{ Row = 51; B = "<which>"; C = 62 },
    @{ Row = 52; B = "<an>"; C = 54 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}
